$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in Quantity column (E) for rows 10 and 11 that were previously blank ---
$ws.Range("E10").Value = "NA"
$ws.Range("E11").Value = 1

# --- Row 12 : Item 11 ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Black-Oxide 18-8 Pan Head Phillips Screws, M2.5 x 0.45 mm Thread, 10mm Long"
$ws.Range("C12").Value = "McMaster Carr"
$ws.Range("D12").Value = "95836A257"
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "https://www.mcmaster.com/95836A257/?SrchEntryWebPart_InpBox=95836A257"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"

# --- Row 13 : Item 12 ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Steel Thin Hex Nut Medium-Strength, M2.5 x 0.45 mm "
$ws.Range("C13").Value = "McMaster Carr"
$ws.Range("D13").Value = "90370A202"
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = "https://www.mcmaster.com/90370A202/?SrchEntryWebPart_InpBox=95836A257"
$ws.Range("G13").Value = "NA"
$ws.Range("H13").Value = "NA"

# --- Row 14 : Item 13 ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Phillips Screws M3 x 0.50 mm Thread, 20mm Long"
$ws.Range("C14").Value = "McMaster Carr"
$ws.Range("D14").Value = "95836A535"
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = "https://www.mcmaster.com/95836A535/?SrchEntryWebPart_InpBox=95836A257"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "NA"

# --- Row 15 : Item 14 (B15 already held the "not quite done yet" placeholder text) ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Hex Nut M3 x 0.5 mm Thread"
$ws.Range("C15").Value = "McMaster Carr"
$ws.Range("D15").Value = "98676A100"
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "https://www.mcmaster.com/98676A100/?SrchEntryWebPart_InpBox=95836A257"
$ws.Range("G15").Value = "NA"
$ws.Range("H15").Value = "NA"

# --- Row 16 : Item 15 ---
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Black-Oxide 18-8 Pan Head Phillips Screws M2.5 x 0.45 mm Thread, 20 mm Long"
$ws.Range("C16").Value = "McMaster Carr"
$ws.Range("D16").Value = "95836A216"
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = "https://www.mcmaster.com/95836A216/?SrchEntryWebPart_InpBox=95836A257"
$ws.Range("G16").Value = "NA"
$ws.Range("H16").Value = "NA"

# --- Widen columns B and F to fit the new, longer descriptions/URLs ---
$ws.Columns("B").ColumnWidth = 68.8
$ws.Columns("F").ColumnWidth = 67.1

# --- Move the active selection to reflect where editing left off ---
$ws.Range("B26").Select() | Out-Null
